$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: restructure rows, preserving the special styled rows (1, 6->5, 17->15, 18->16) ---
# Delete rows bottom-up that carry no special formatting (their values will be rewritten below).
$ws.Rows(35).Delete()
$ws.Rows(33).Delete()
$ws.Rows(31).Delete()
$ws.Rows(29).Delete()
$ws.Rows(27).Delete()
$ws.Rows(25).Delete()
$ws.Rows(23).Delete()
$ws.Rows(21).Delete()
$ws.Rows(19).Delete()
$ws.Rows(12).Delete()
$ws.Rows(10).Delete()
$ws.Rows(8).Delete()
$ws.Rows(4).Delete()
$ws.Rows(2).Delete()

# Re-open the gaps needed so the remaining "special" rows land on their final target rows.
$ws.Rows(4).Insert()
$ws.Rows(13).Insert()
$ws.Rows(14).Insert()

# --- Step 2: write the new cell content into place ---

# Row 1 (unchanged): Parameters / Control Group / Treatment Group
$ws.Range("A1").Value = "Parameters"
$ws.Range("B1").Value = "Control Group"
$ws.Range("C1").Value = "Treatment Group"

# Row 3: Sample Size (n)
$ws.Range("A3").Value = "Sample Size (n)"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1.5

# Row 5: Population Mean / Total Number of Simulations (keeps ht=29 from original row 6)
$ws.Range("A5").Value = "Population Mean"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("F5").Value = "Total Number of Simulations:"
$ws.Range("G5").Value = 10000

# Row 7: Population Variance
$ws.Range("A7").Value = "Population Variance"
$ws.Range("B7").Value = 30
$ws.Range("C7").Value = 30

# Row 9: Population Distribution
$ws.Range("A9").Value = "Population Distribution"
$ws.Range("B9").Value = "Normal (Standard)"
$ws.Range("C9").Value = "Normal"

# Row 11: Sample Variance
$ws.Range("A11").Value = "Sample Variance"
$ws.Range("B11").Value = "sX^2"
$ws.Range("C11").Value = "sY^2"

# Row 13: Sample Mean
$ws.Range("A13").Value = "Sample Mean"
$ws.Range("B13").Value = "x_bar"
$ws.Range("C13").Value = "y_bar"

# Row 15: Statistical Test header (kept s=1 customFormat from original row 17)
$ws.Range("A15").Value = "Statistical Test"
$ws.Range("B15").Value = "F Ratio Test"
$ws.Range("C15").Value = "Levene's Test"

# Row 16 stays blank (kept s=3 customFormat from original row 18)

# Row 17: Significance Level (alpha)
$ws.Range("A17").Value = "Significance Level (alpha)"
$ws.Range("B17").Value = 0.05
$ws.Range("C17").Value = 0.05

# Row 19: Test Statistic
$ws.Range("A19").Value = "Test Statistic "
$ws.Range("B19").Value = "F = sX^2/sY^2"

# Row 21: Null Distribution (H0)
$ws.Range("A21").Value = "Null Distribution (H0)"
$ws.Range("B21").Value = "F_30-1,30-1"

# Row 23: Alternative Distribution (Ha)
$ws.Range("A23").Value = "Alternative Distribution (Ha)"
$ws.Range("B23").Value = "F_30-1,30-1(ncp = 0+0)"

# Row 25: Critical Value (lower)
$ws.Range("A25").Value = "Critical Value"
$ws.Range("B25").Value = "F_0.05, 30-1, 30-1 (Lower)"
$ws.Rows(25).RowHeight = 29

# Row 26: Critical Value (upper)
$ws.Range("B26").Value = "F_1-0.05, 30-1, 30-1 (Upper)"
$ws.Rows(26).RowHeight = 29

# Row 27: Critical Value (two-sided)
$ws.Range("B27").Value = "F_1-0.05/2, 30-1, 30-1 and F_0.05/2, 30-1, 30-1 (Two)"
$ws.Rows(27).RowHeight = 43.5

# Row 29: Type I Error Rate
$ws.Range("A29").Value = "Type I Error Rate"
$ws.Range("B29").Value = 0.05

# Row 31: Type II Error Rate
$ws.Range("A31").Value = "Type II Error Rate"

# Row 33: Theoretical Power
$ws.Range("A33").Value = "Theoretical Power"

# Row 35: Empirical Power
$ws.Range("A35").Value = "Empirical Power"
$ws.Range("B35").Formula = "=1233/10000"

# --- Step 3: sheet view selection ---
$ws.Range("B33").Select()

Write-Output "done"
